# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 64 (pushing the existing rows 64-89
# down to 65-90) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 64..89 down to 65..90, creating a fresh (formatted) row 64.
$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value = "2022-01-06"
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112008
$ws.Cells.Item(64, 7).Value = "Coliflor"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Tercera"
$ws.Cells.Item(64, 10).Value = 1200
$ws.Cells.Item(64, 11).Value = 400
$ws.Cells.Item(64, 12).Value = 500
$ws.Cells.Item(64, 13).Value = 450
$ws.Cells.Item(64, 14).Value = "$/unidad"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 450
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
